# Full manual control finished
# that took longer than i expected...
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# PWM labels: GearRail now reports Invalid, GearGate now reports PWM-6
$ws.Range("G5").Value = "GearRail CAN-12 PWM-Invalid"
$ws.Range("G6").Value = "GearGate CAN-7 PWM-6"

# Joystick 1 button B now triggers PlaceGear; Joystick 2 actions renamed
$ws.Range("B13").Value = "PlaceGear"
$ws.Range("D13").Value = "Gear Gate Close"
$ws.Range("D14").Value = ""

# Aim marked as not-yet-implemented; Start button now opens the gear chute
$ws.Range("B21").Value = "Aim (still to do)"
$ws.Range("D21").Value = "GearChute"

# Restore the editor's last selection
$ws.Range("E19").Select()
